# Rename the "String" column header to "String(s)" to reflect that the
# column now supports the new non-strategic list types (GRC, WPAFCC,
# UT_IAC), and restore the active selection to E6 as left by the author.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = "String(s)"

$ws.Range("E6").Select()
